$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "rating" column (old C),
# shifting it to column E and leaving blank C/D in between.
$ws.Columns("C:D").Insert()

# --- Header row (row 1) ---
# B1 keeps "UN-header" style column, gets the newest date label.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"
$ws.Range("E1").Value = "Jun_10"

# --- Data rows (2-27): new C & D columns mirror column B's "UN" values ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- Column widths: keep the same custom width (8) across C, D, E ---
$ws.Columns("C").ColumnWidth = 7.140625
$ws.Columns("D").ColumnWidth = 7.140625
$ws.Columns("E").ColumnWidth = 7.140625
